$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:A4").ClearContents()
$ws.Range("D2:D4").ClearContents()
$ws.Hyperlinks.Delete()
